$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("H6").Select()
